$d = $word.ActiveDocument
$sec = $d.Sections(1)

function Rename-InlineImage {
    param($Range, $DocPrId, $CNvPrId, $Descr, $NewName, $ExtentCx, $ExtentCy, $EmbedRId)

    $xmlFrag = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:r><w:drawing><wp:inline distB="0" distT="0" distL="0" distR="0"><wp:extent cx="$ExtentCx" cy="$ExtentCy"/><wp:effectExtent b="0" l="0" r="0" t="0"/><wp:docPr descr="$Descr" id="$DocPrId" name="$NewName"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="$Descr" id="$CNvPrId" name="$NewName"/><pic:cNvPicPr preferRelativeResize="0"/></pic:nvPicPr><pic:blipFill><a:blip r:embed="$EmbedRId"/><a:srcRect b="0" l="0" r="0" t="0"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="$ExtentCx" cy="$ExtentCy"/></a:xfrm><a:prstGeom prst="rect"/><a:ln/></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@

    $shp = $Range.InlineShapes(1)
    $r = $shp.Range
    $shp.Delete()
    $r.InsertXML($xmlFrag)
}

# Footer, first page (footer1.xml): PearsonLogo id="3"/"0", image1.png -> image2.png
$footerFirst = $sec.Footers(2)
Rename-InlineImage $footerFirst.Range "3" "0" 'Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png' "image2.png" "952500" "285750" "rId1"

# Footer, default (footer2.xml): PearsonLogo id="2"/"0", image1.png -> image2.png
$footerDefault = $sec.Footers(1)
Rename-InlineImage $footerDefault.Range "2" "0" 'Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png' "image2.png" "952500" "285750" "rId1"

# Header, first page (header1.xml): BTec_Logo-Orange id="1"/"0", image2.jpg -> image1.jpg
$headerFirst = $sec.Headers(2)
Rename-InlineImage $headerFirst.Range "1" "0" "BTec_Logo-Orange" "image1.jpg" "914400" "277792" "rId1"
